# "Add files via upload" — current versions of structural tables, to be used
# during "non-presentation".
#
# The data that used to live in column M ("comms") for every data row is
# moved over to column S ("comms_internal"); column M is left blank. The
# cell formatting follows the value, so every M/S cell touched ends up
# sharing the same look (style index 5 in the original file).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 131
$colM = 13   # column M — comms
$colS = 19   # column S — comms_internal

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $mCell = $ws.Cells.Item($r, $colM)
    $sCell = $ws.Cells.Item($r, $colS)

    # Move the value from M into S (S was empty before the edit).
    $sCell.Value = $mCell.Value()

    # Make S pick up M's formatting (so both land on the same style),
    # then clear M's value while keeping its own formatting untouched.
    $mCell.Copy()
    $sCell.PasteSpecial(-4122)   # xlPasteFormats
    $mCell.ClearContents()
}

$excel.CutCopyMode = 0

# Reflect the new area of interest in the sheet view: the cursor now lands
# on M2 (top of the now-empty column) with the whole M2:M131 range marked
# as selected, and the window is scrolled so column H is the first visible
# column.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 8
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("M2:M131").Select()
